# Update Data Sources from LFX
#
# The refreshed data-source tables ship with a new PowerPoint table style
# (GUID {928D4016-75EB-4100-83A0-C58B1C6B50E9}) in place of the previous
# default style (GUID {8FA931E9-57A9-4F2B-8CB2-13A407556EF8}). Walk every
# slide/shape in the deck and re-apply the new style to every table we find.

$oldStyleId = "{8FA931E9-57A9-4F2B-8CB2-13A407556EF8}"
$newStyleId = "{928D4016-75EB-4100-83A0-C58B1C6B50E9}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)

        if ($shape.HasTable) {
            $table = $shape.Table

            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
